$d = $word.ActiveDocument

# The development diary table is the first (only) table in the document.
$t = $d.Tables.Item(1)

# Append a new row, copying formatting (shading/borders/widths) from the
# last existing row, then populate its cells with the new entry.
$newRow = $t.Rows.Add()

$newRow.Cells.Item(1).Range.Text = "23/11/2021"
$newRow.Cells.Item(2).Range.Text = "2 Hours"
$newRow.Cells.Item(3).Range.Text = "Stage 3"
$newRow.Cells.Item(4).Range.Text = "Implemented the input recognition system and tested the usage of beginpaint() and the HDC system. Some additions are still required to finish the stage however the basic foundations for the procedure are present."
